$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new ("fake combined") row of tyre data below the existing table.
# H/I are set before A so the new shared strings register in the same
# order the author's session produced them in.
$ws.Range("H30").Value = "1,2"
$ws.Range("I30").Value = "3,4"
$ws.Range("A30").Value = "'29"
$ws.Range("B30").Value = "Hoosier"
$ws.Range("C30").Value = "16.0x6.0-10"
$ws.Range("D30").Value = 16
$ws.Range("E30").Value = "LCO"
$ws.Range("F30").Value = 6
$ws.Range("G30").Value = 99
$ws.Range("J30").Formula = '=_xlfn.CONCAT( A30, "_", B30,"_",C30,"_",E30,"_",F30,"Rim.tir")'
$ws.Range("K30").Value = 5.1323934799999993

# Match the formatting used by the rest of the table (row 16 has the same
# per-column style pattern the new row needs for columns B:K).
$ws.Range("B16:K16").Copy()
$ws.Range("B30:K30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the view state to match where the author was looking when saving.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("K24").Select()
